$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cell values as per diff.
# Column D values are numeric-looking strings (e.g. "1.001", "28.086.49")
# stored as text in the source data, so force text format per-cell before
# assigning to prevent Excel auto-converting them to numbers.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.086.49'
$ws.Range('E2').Value = '  -2.11%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.834.92'
$ws.Range('E3').Value = '  -0.94%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '324.62'
$ws.Range('E5').Value = '  -2.99%  '
$ws.Range('E6').Value = '  -0.12%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4635'
$ws.Range('E7').Value = '  -0.46%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3868'
$ws.Range('E8').Value = '  -1.22%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07863'
$ws.Range('E10').Value = '  -2.49%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '21.94'
$ws.Range('E11').Value = '  -1.91%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.804.46'
$ws.Range('E12').Value = '  -5.13%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.688'
$ws.Range('E13').Value = '  -2.79%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.919'
$ws.Range('E14').Value = '  -1.19%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06847'
$ws.Range('E15').Value = '  +0.02%  '
$ws.Range('E16').Value = '  -0.53%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.001'
$ws.Range('E17').Value = '  -0.10%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009956'
$ws.Range('E18').Value = '  -1.47%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.70'
$ws.Range('E19').Value = '  -2.67%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.000'
$ws.Range('E20').Value = '  -0.18%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '28.106.08'
$ws.Range('E21').Value = '  -2.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.331'
$ws.Range('E22').Value = '  -1.18%  '
$ws.Range('E23').Value = '  -2.85%  '
$ws.Range('E24').Value = '  -2.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.006.22'
$ws.Range('E25').Value = '  -5.02%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '154.01'
$ws.Range('E26').Value = '  +0.49%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.18'
$ws.Range('E27').Value = '  -1.40%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.696'
$ws.Range('E28').Value = '  -6.96%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.966'
$ws.Range('E29').Value = '  -2.98%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '118.21'
$ws.Range('E30').Value = '  +0.40%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.9379'
$ws.Range('E31').Value = '  -3.99%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09247'
$ws.Range('E32').Value = '  -2.10%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.273'
$ws.Range('E33').Value = '  -1.89%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.324'
$ws.Range('E34').Value = '  -2.27%  '
$ws.Range('E35').Value = '  -5.24%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.05856'
$ws.Range('E36').Value = '  -5.20%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02129'
$ws.Range('E37').Value = '  -3.42%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.142'
$ws.Range('E38').Value = '  -1.93%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '7.777'
$ws.Range('E39').Value = '  +2.11%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5597'
$ws.Range('E40').Value = '  -2.31%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '9.896'
$ws.Range('E41').Value = '  -2.82%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1762'
$ws.Range('E42').Value = '  -2.05%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.07231'
$ws.Range('E43').Value = '  +1.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '11.61'
$ws.Range('E44').Value = '  -1.97%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5269'
$ws.Range('E45').Value = '  -2.63%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.136'
$ws.Range('E46').Value = '  -10.08%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.120'
$ws.Range('E47').Value = '  -10.43%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.835'
$ws.Range('E48').Value = '  -4.26%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '112.62'
$ws.Range('E49').Value = '  -1.38%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.9992'
$ws.Range('E50').Value = '  -0.10%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.323'
$ws.Range('E51').Value = '  +0.34%  '
